$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The "COMPRADOR" column (column 3) values are rotating between rows.
# Set each target cell directly by row/column index to avoid any
# ambiguity/collision that a text-based Find/Replace would have when
# several cells share swapped values.

$t.Cell(2, 3).Range.Text  = "INTERLINK2AMERICAS"   # was FLOR A FRUTO
$t.Cell(3, 3).Range.Text  = "COLFRESH COFFEE"      # was CAFÉ MOLINA
$t.Cell(4, 3).Range.Text  = "CAFÉ MOLINA"          # was REGIONAL S.A.S
$t.Cell(5, 3).Range.Text  = "REGIONAL S.A.S"       # was ARMANDO VELÁSQUEZ
$t.Cell(6, 3).Range.Text  = "FLOR A FRUTO"         # was COLFRESH COFFEE
$t.Cell(7, 3).Range.Text  = "INMERSSO BOUTIQUE"    # was BOX BRAND
$t.Cell(8, 3).Range.Text  = "ARMANDO VELÁSQUEZ"    # was NEIRA YORK COFFEE
$t.Cell(10, 3).Range.Text = "BOX BRAND"            # was INTERLINK2AMERICAS
$t.Cell(11, 1).Range.Text = "11:30 - 11:45"        # was 11:45 - 12:00
$t.Cell(11, 3).Range.Text = "NEIRA YORK COFFEE"    # was INMERSSO BOUTIQUE
